$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the text values "Patidev" / "Rishabh" with the number 123
$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 123
